$d = $word.ActiveDocument

# --- 1. Fix typo in heading: RESTRINCCIONES -> RESTRICCIONES ---
$d.Content.Find.Execute(".1. RESTRINCCIONES DEL SOFTWARE", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ".1. RESTRICCIONES DEL SOFTWARE", 2)

# --- 2. Accent fixes on existing bullet paragraphs ---
$d.Content.Find.Execute(
    "Se validaran las relaciones de las tablas para prohibir agregar codigos inexistentes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Se validaran las relaciones de las tablas para prohibir agregar códigos inexistentes", 2)

$d.Content.Find.Execute(
    "Los codigos casi en su totalidad seran solamente valores numericos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Los códigos casi en su totalidad serán solamente valores numéricos", 2)

$d.Content.Find.Execute(
    "No se permitira borrar ningun dato ya almacenado, solo se podria inactivar o cancelar.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "No se permitirá borrar ningún dato ya almacenado, solo se podría inactivar o cancelar.", 2)

# --- 3. Append three new bullet paragraphs after the last ("Los datos
#        presentados...") paragraph, and move the "_GoBack" bookmark out
#        of that paragraph into the new final one, splitting its text
#        across two runs around the bookmark ("...mismo e-" | "mail
#        registrado."), matching the authored edit. ---

# The "_GoBack" bookmark currently sits at the start of the last
# paragraph; drop it here so it can be re-created at its new home
# (a document can only usefully have one "_GoBack").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
# Collapsed insertion point immediately before the paragraph's own
# end-of-paragraph mark, so the new paragraphs land right after it
# without leaving any stray empty paragraph behind.
$insertPoint = $d.Range($r.End - 1, $r.End - 1)

$xmlFragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-DO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>Los usuarios deben estar registrados si desean reemplazar productos de receta o agregar productos al carro de compras.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-DO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>Cada receta y producto debe tener la foto para ser presentados al usuario</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-DO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>No puede existir mas de un usuario con el mismo e-</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>mail registrado.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xmlFragment)

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
